# Generate Report for Handoff
# Updates the localization-status report: new source file GUID
# (bf14b830-...) replaces the previous run's GUID (e7a4e73a-...),
# refreshed handoff timestamps/xlf names, and the now-empty
# "Latest Target File" / "Latest Handback File" columns (handback
# hasn't happened yet for this run) on the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
$wsOverview.Range("B2").Value = "e2e\bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
foreach ($h in $wsOverview.Hyperlinks) {
    if ($h.Range.Address() -eq '$B$2') {
        $h.TextToDisplay = "e2e\bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
    }
}
$wsOverview.Range("G2").Value = "2016-08-28 08:58:30"

# ---------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
    }
}
$wsZhCn.Range("G2").Value = "bf14b830-d9bf-40d2-9484-9712823cfb6b.b8900584d9d9a94e7d00b14d6e4b96741c82913a.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 08:58:26"

# Latest Target File hyperlink/value is cleared - no handback yet for this run
foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$I$2') {
        $h.Delete()
    }
}
$wsZhCn.Range("I2").ClearContents()
$wsZhCn.Range("I2").Style = "Normal"

# Latest Handback File is cleared too
$wsZhCn.Range("J2").ClearContents()

# Latest Handback DateTime resets to the zero-date sentinel
$wsZhCn.Range("K2").Value = "0001-01-01 00:00:00"

$wsZhCn.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsZhCn.Columns.Item(10).ColumnWidth = 21.7054770333426

# ---------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "bf14b830-d9bf-40d2-9484-9712823cfb6b.md"
    }
}
$wsDeDe.Range("G2").Value = "bf14b830-d9bf-40d2-9484-9712823cfb6b.b8900584d9d9a94e7d00b14d6e4b96741c82913a.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 08:58:30"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$I$2') {
        $h.Delete()
    }
}
$wsDeDe.Range("I2").ClearContents()
$wsDeDe.Range("I2").Style = "Normal"

$wsDeDe.Range("J2").ClearContents()

$wsDeDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDeDe.Columns.Item(9).ColumnWidth = 18.6506053379604
$wsDeDe.Columns.Item(10).ColumnWidth = 21.7054770333426
